$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Remove the "Percent of variance by level" row entirely - everything below
# shifts up one row (D4's lone numeric value disappears with it; it is
# re-introduced, shifted, as D3 below).
$ws.Rows("4:4").Delete()

# Row 3: renamed section label, plus a new numeric flag in column D.
$ws.Range("A3").Value = "Baseline Variance"
$ws.Range("D3").Value = 1

# Row 2 header labels gain sample-size annotations / get simplified.
$ws.Range("B2").Value = "Team-Levela (n = 87)"
$ws.Range("C2").Value = "Model-Levela (n = 1,253)"
$ws.Range("D2").Value = "Total"

# Row 13 becomes a new sub-heading bullet ("Belief in Hypothesis"), styled
# like the other section headers (left-aligned) instead of a plain item.
$ws.Range("A13").Value = "Belief in Hypothesis"
$ws.Range("A13").HorizontalAlignment = -4131

# Rows 14-17 get reordered / relabeled.
$ws.Range("A14").Value = "Pro-Immigration Attitude"
$ws.Range("A15").Value = "Knowledge of Topic"
$ws.Range("A16").Value = "Statistical Skills"
$ws.Range("A17").Value = "Peer Model Ranking"

# "Unexplained Variance" moves down one row, leaving row 18 blank again.
$ws.Range("A18").Value = ""
$ws.Range("A19").Value = "Unexplained Variance"
$ws.Range("A19").HorizontalAlignment = -4152

# Sheet3 becomes the active tab/sheet, with B19 selected; Sheet1 loses focus.
$ws.Activate()
$ws.Range("B19").Select()
